# "major model revamp with pickups"
# Add a weekly pickup-frequency column of data under the existing
# (renamed) "Fréquence de Ramasse(/w)" header, and widen the column to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header of column E to indicate the frequency is per week.
$ws.Range("E1").Value = "Fréquence de Ramasse(/w)"

# Fill in the weekly pickup frequency for each pickup point.
$ws.Range("E2").Value = 4
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("E8").Value = 2

# These are whole numbers of pickups per week - format as integers.
$ws.Range("E2:E8").NumberFormat = "#,##0"

# Make sure the wider header text is fully visible.
$ws.Columns.Item(5).EntireColumn.AutoFit()
